$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2022-09-15 21:02:36"

for ($row = 2; $row -le 63; $row++) {
    $ws.Range("O$row").Value = $newTimestamp
}
